# Daily attendance processing - 2025-11-04 21:43:22
# Normalizes the "Recorded By" (column G) cell values so that any
# "System"/"system" entries are moved to the front of the comma-separated
# list, preserving the relative order of the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($systemParts.Count -eq 0) {
        continue
    }

    $newParts = $systemParts + $otherParts
    $newVal = $newParts -join ", "

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
